$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Browser"/"Chrome" column (column B), shifting the
# remaining columns (Search/Search_2, Apple) one to the left.
$ws.Range("B1").EntireColumn.Delete()

$ws.Range("B1").Select()
